$wb = $excel.ActiveWorkbook

# --- 1. Update the Date value on the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- 2. Swap columns AK and AL on the Elements sheet ---
$ws = $wb.Worksheets.Item("Elements")

# Save current values for column AK (37) and AL (38), rows 1-6
$akValues = @{}
$alValues = @{}
for ($r = 1; $r -le 6; $r++) {
    $akValues[$r] = $ws.Cells.Item($r, 37).Value()
    $alValues[$r] = $ws.Cells.Item($r, 38).Value()
}

# Write swapped values back (AK gets AL's old content, AL gets AK's old content)
# Only touch cells whose value actually changes, to avoid disturbing rows
# where AK and AL already held identical (e.g. both empty) content.
for ($r = 1; $r -le 6; $r++) {
    if ($akValues[$r] -ne $alValues[$r]) {
        $ws.Cells.Item($r, 37).Value = $alValues[$r]
        $ws.Cells.Item($r, 38).Value = $akValues[$r]
    }
}

# Swap column widths (AK was 24.98046875, AL was 64.890625 -> now swapped)
$ws.Columns.Item(37).ColumnWidth = 64.0
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668
